# Poll Coverage Tooba.xlsx - check-in update
#
# 1) B40 had a blank/zero data point for poll #38; fill in the observed value.
# 2) Select the cell that was just edited (mirrors the authored selection move).
# 3) Recolor the workbook's theme "Background 1" (lt1) swatch to the light
#    green used in the refreshed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit -------------------------------------------------------
$ws.Range("B40").Value = 0.63

# --- Selection / navigation -------------------------------------------
$ws.Range("B40").Select() | Out-Null

# --- Theme color update (lt1 / Background 1 -> 78DC78) ---------------
$themeColors = $wb.Theme.ThemeColorScheme
$background1 = $themeColors.Colors(2)
$background1.RGB = 0x78DC78
